# Fix the misspelled "Mananger" -> "Manager" in the three
# "Resource Manager" rounded-rectangle shapes on slide 12.
# The original text is split across two runs ("Resource " + "Mananger");
# we rewrite the whole text range so PowerPoint merges it back into a
# single run (using the first run's formatting) reading "Resource Manager".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)

$shapeIndexes = @(10, 16, 18)

foreach ($idx in $shapeIndexes) {
    $shp = $s.Shapes.Item($idx)
    $tr = $shp.TextFrame.TextRange
    $len = $tr.Characters().Count
    $full = $tr.Characters(1, $len)
    $full.Text = "Resource Manager"
}
